# Applies the scheduled-runner profit/price refresh to the Titan_Profits workbook.
# Generated from the canonical OOXML diff: per-row numeric refresh across the
# ALC / ARM / BSM / CRP / GSM / LTW / WVR leve-profit tables (CUL unchanged).

$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
# row 31
$ws.Range("H31").Value = 3503
$ws.Range("I31").Value = 3503
$ws.Range("K31").Value = 10509
$ws.Range("M31").Value = -10279
# row 38
$ws.Range("H38").Value = 324.81818
$ws.Range("I38").Value = 60.5
$ws.Range("K38").Value = 181.5
$ws.Range("M38").Value = 190.5
# row 40
$ws.Range("H40").Value = 2242.2856
$ws.Range("I40").Value = 2340
$ws.Range("J40").Value = 2144.5715
$ws.Range("K40").Value = 2340
$ws.Range("L40").Value = 2144.5715
$ws.Range("M40").Value = -2165
$ws.Range("N40").Value = -2494.5715
# row 55
$ws.Range("H55").Value = 156.41667
$ws.Range("I55").Value = 109.375
$ws.Range("J55").Value = 250.5
$ws.Range("K55").Value = 109.375
$ws.Range("L55").Value = 250.5
$ws.Range("M55").Value = 104.625
$ws.Range("N55").Value = -678.5
# row 64
$ws.Range("H64").Value = 3879.6667
$ws.Range("I64").Value = 3622.5
$ws.Range("J64").Value = 4222.5557
$ws.Range("K64").Value = 3622.5
$ws.Range("L64").Value = 4222.5557
$ws.Range("M64").Value = -3374.5
$ws.Range("N64").Value = -4718.5557
# row 67
$ws.Range("H67").Value = 3879.6667
$ws.Range("I67").Value = 3622.5
$ws.Range("J67").Value = 4222.5557
$ws.Range("K67").Value = 3622.5
$ws.Range("L67").Value = 4222.5557
$ws.Range("M67").Value = -2764.5
$ws.Range("N67").Value = -5938.5557
# row 68
$ws.Range("H68").Value = 17500
$ws.Range("J68").Value = 17500
$ws.Range("L68").Value = 17500
$ws.Range("N68").Value = -18998
# row 71
$ws.Range("H71").Value = 17500
$ws.Range("J71").Value = 17500
$ws.Range("L71").Value = 52500
$ws.Range("N71").Value = -59988
# row 74
$ws.Range("H74").Value = 3625.5833
$ws.Range("I74").Value = 3560.3
$ws.Range("J74").Value = 3952
$ws.Range("K74").Value = 3560.3
$ws.Range("L74").Value = 3952
$ws.Range("M74").Value = -2624.3
$ws.Range("N74").Value = -5824
# row 76
$ws.Range("H76").Value = 3706987
$ws.Range("I76").Value = 4632742
$ws.Range("J76").Value = 3966.6667
$ws.Range("K76").Value = 4632742
$ws.Range("L76").Value = 3966.6667
$ws.Range("M76").Value = -4632427
$ws.Range("N76").Value = -4596.6667
# row 77
$ws.Range("H77").Value = 3625.5833
$ws.Range("I77").Value = 3560.3
$ws.Range("J77").Value = 3952
$ws.Range("K77").Value = 17801.5
$ws.Range("L77").Value = 19760
$ws.Range("M77").Value = -13121.5
$ws.Range("N77").Value = -29120
# row 79
$ws.Range("H79").Value = 3706987
$ws.Range("I79").Value = 4632742
$ws.Range("J79").Value = 3966.6667
$ws.Range("K79").Value = 4632742
$ws.Range("L79").Value = 3966.6667
$ws.Range("M79").Value = -4631650
$ws.Range("N79").Value = -6150.6667
# row 92
$ws.Range("H92").Value = 856499.3
$ws.Range("I92").Value = 1236223.2
$ws.Range("J92").Value = 2120.5
$ws.Range("K92").Value = 1236223.2
$ws.Range("L92").Value = 2120.5
$ws.Range("M92").Value = -1234975.2
$ws.Range("N92").Value = -4616.5
# row 95
$ws.Range("H95").Value = 20000
$ws.Range("J95").Value = 20000
$ws.Range("L95").Value = 20000
$ws.Range("N95").Value = -25492
# row 97
$ws.Range("H97").Value = 1900
$ws.Range("J97").Value = 1900
$ws.Range("L97").Value = 5700
$ws.Range("N97").Value = -6692
# row 138
$ws.Range("H138").Value = 10896229
$ws.Range("I138").Value = 2780847.2
$ws.Range("J138").Value = 16131960
$ws.Range("K138").Value = 8342541.600000001
$ws.Range("L138").Value = 48395880
$ws.Range("M138").Value = -8337401.600000001
$ws.Range("N138").Value = -48406160

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
# row 63
$ws.Range("H63").Value = 4233.0435
$ws.Range("I63").Value = 4097.9473
$ws.Range("J63").Value = 4874.75
$ws.Range("K63").Value = 4097.9473
$ws.Range("L63").Value = 4874.75
$ws.Range("M63").Value = -3411.9473
$ws.Range("N63").Value = -6246.75
# row 66
$ws.Range("H66").Value = 4233.0435
$ws.Range("I66").Value = 4097.9473
$ws.Range("J66").Value = 4874.75
$ws.Range("K66").Value = 20489.7365
$ws.Range("L66").Value = 24373.75
$ws.Range("M66").Value = -17057.7365
$ws.Range("N66").Value = -31237.75
# row 132
$ws.Range("H132").Value = 3141.5557
$ws.Range("I132").Value = 2663.5173
$ws.Range("K132").Value = 7990.5519
$ws.Range("M132").Value = -5460.5519
# row 133
$ws.Range("H133").Value = 43476.668
$ws.Range("J133").Value = 43476.668
$ws.Range("L133").Value = 43476.668
$ws.Range("N133").Value = -48536.668

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
# row 86
$ws.Range("H86").Value = 1459.65
$ws.Range("I86").Value = 1548.4445
$ws.Range("J86").Value = 1387
$ws.Range("K86").Value = 1548.4445
$ws.Range("L86").Value = 1387
$ws.Range("M86").Value = -425.4445000000001
$ws.Range("N86").Value = -3633
# row 89
$ws.Range("H89").Value = 1459.65
$ws.Range("I89").Value = 1548.4445
$ws.Range("J89").Value = 1387
$ws.Range("K89").Value = 7742.2225
$ws.Range("L89").Value = 6935
$ws.Range("M89").Value = -2126.2225
$ws.Range("N89").Value = -18167
# row 94
$ws.Range("H94").Value = 627.12
$ws.Range("I94").Value = 507.15
$ws.Range("K94").Value = 507.15
$ws.Range("M94").Value = -56.14999999999998
# row 105
$ws.Range("H105").Value = 280892.56
$ws.Range("I105").Value = 3000
$ws.Range("J105").Value = 912466.5600000001
$ws.Range("K105").Value = 3000
$ws.Range("L105").Value = 912466.5600000001
$ws.Range("M105").Value = -1253
$ws.Range("N105").Value = -915960.5600000001

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
# row 58
$ws.Range("H58").Value = 3164.95
$ws.Range("I58").Value = 1212.4445
$ws.Range("J58").Value = 4762.4546
$ws.Range("K58").Value = 1212.4445
$ws.Range("L58").Value = 4762.4546
$ws.Range("M58").Value = -1009.4445
$ws.Range("N58").Value = -5168.4546
# row 62
$ws.Range("H62").Value = 31556.572
$ws.Range("I62").Value = 68033.336
$ws.Range("J62").Value = 4199
$ws.Range("K62").Value = 68033.336
$ws.Range("L62").Value = 4199
$ws.Range("M62").Value = -67409.336
$ws.Range("N62").Value = -5447
# row 65
$ws.Range("H65").Value = 31556.572
$ws.Range("I65").Value = 68033.336
$ws.Range("J65").Value = 4199
$ws.Range("K65").Value = 340166.68
$ws.Range("L65").Value = 20995
$ws.Range("M65").Value = -337046.68
$ws.Range("N65").Value = -27235
# row 136
$ws.Range("H136").Value = 3164.95
$ws.Range("I136").Value = 1212.4445
$ws.Range("J136").Value = 4762.4546
$ws.Range("K136").Value = 3637.3335
$ws.Range("L136").Value = 14287.3638
$ws.Range("M136").Value = -1087.3335
$ws.Range("N136").Value = -19387.3638

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
# row 80
$ws.Range("H80").Value = 2925
$ws.Range("I80").Value = 2957.1428
$ws.Range("J80").Value = 2850
$ws.Range("K80").Value = 2957.1428
$ws.Range("L80").Value = 2850
$ws.Range("M80").Value = -1959.1428
$ws.Range("N80").Value = -4846
# row 83
$ws.Range("H83").Value = 2925
$ws.Range("I83").Value = 2957.1428
$ws.Range("J83").Value = 2850
$ws.Range("K83").Value = 14785.714
$ws.Range("L83").Value = 14250
$ws.Range("M83").Value = -9793.714
$ws.Range("N83").Value = -24234
# row 132
$ws.Range("H132").Value = 3318.5
$ws.Range("I132").Value = 3430.8572
$ws.Range("K132").Value = 10292.5716
$ws.Range("M132").Value = -7762.571599999999

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
# row 55
$ws.Range("H55").Value = 436.53333
$ws.Range("J55").Value = 452.54544
$ws.Range("L55").Value = 452.54544
$ws.Range("N55").Value = -798.54544
# row 132
$ws.Range("H132").Value = 5369.136
$ws.Range("I132").Value = 3300.3333
$ws.Range("K132").Value = 9900.999899999999
$ws.Range("M132").Value = -7370.999899999999
# row 133
$ws.Range("H133").Value = 50487
$ws.Range("J133").Value = 50487
$ws.Range("L133").Value = 50487
$ws.Range("N133").Value = -55547
# row 137
$ws.Range("H137").Value = 39600
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 39600
$ws.Range("K137").Value = 0
$ws.Range("M137").ClearContents()
$ws.Range("N137").Value = -49800

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
# row 64
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("N64").ClearContents()
$ws.Range("L64").Value = 0
# row 67
$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("N67").ClearContents()
$ws.Range("L67").Value = 0
# row 70
$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("M70").ClearContents()
# row 73
$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("M73").ClearContents()
# row 76
$ws.Range("H76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("N76").ClearContents()
$ws.Range("L76").Value = 0
# row 79
$ws.Range("H79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("N79").ClearContents()
$ws.Range("L79").Value = 0
# row 132
$ws.Range("H132").Value = 3515.3684
$ws.Range("I132").Value = 3386.2
$ws.Range("K132").Value = 10158.6
$ws.Range("M132").Value = -7628.599999999999

